# Update "Förändrad" (changed) date column C for rows 2-11 from 45175 (2023-09-06)
# to 45183 (2023-09-14), matching the automatic update of the logging file.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45183
}
